# cht-conf regenerates this XLSForm and, as part of that refresh, drops the
# stray "NO_LABEL" placeholder that used to sit in the label column of the
# "begin_group my_page" row on the survey sheet (row 3, column C). Clearing
# the cell (rather than just blanking its text) removes it from the sheet
# entirely, which also lets the now-unused "NO_LABEL" shared string drop out
# of the shared-string table and every later string index shift down by one
# — matching how the upstream tool re-exports the workbook.

$wb = $excel.ActiveWorkbook
$survey = $wb.Worksheets.Item("survey")

$survey.Range("C3").Clear()

# The regenerated file also lands on the default top-left cell (A1) instead
# of the previous ad-hoc scroll position (B24) for the bottom-right (frozen)
# pane of the survey sheet.
$survey.Range("A1").Select()
